$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Environment Dev -> Preprod, Password placeholder -> "password", Country Romania -> Australia
# A2/C2 lose their wrap/valign formatting (matching the unformatted style used in row 9),
# so copy that formatting over before setting the new values.
$ws.Range("A9").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "Preprod"

$ws.Range("C9").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "password"

$ws.Range("D2").Value = "Australia"

# Row 9: only Country changes Romania -> Australia
$ws.Range("D9").Value = "Australia"

# Row 10: Password placeholder -> "password", Country Romania -> Australia
# D10 picks up the wrap/valign formatting used elsewhere in column D.
$ws.Range("C10").Value = "password"

$ws.Range("D2").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = "Australia"

# Row 11: only Country changes Romania -> Australia
$ws.Range("D11").Value = "Australia"

$excel.CutCopyMode = $false
